$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 value
$ws.Range("A2").Value = 0.28

# Fill in new row 3 data
$ws.Range("A3").Value = 0.245
$ws.Range("B3").Formula = "=SQRT(1000000/((-0.03069343+SQRT(0.03069343^2-4*0.00074611*(0.20083726-A3)))/(2*0.00074611)))-273.15"

# Update selection
$ws.Range("A3").Select()
